$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -97.341
$ws.Range("B2").Value = -97.3138

$ws.Range("A3").Value = 27.7836
$ws.Range("B3").Value = 27.8079

$ws.Range("A4").Value = -96.7095
$ws.Range("B4").Value = -96.7367

$ws.Range("A5").Value = 28.3675
$ws.Range("B5").Value = 28.3433
